$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.412.38"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "3.145.42"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'561.78"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Value = "'142.12"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.139.94"
$ws.Range("E8").Value = "  +2.91%  "
$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").Value = "'6.76"
$ws.Range("E10").Value = "  +5.41%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "'0.465"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").Value = "'36.61"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "'0.0000221"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "3.645.46"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").Value = "64.488.03"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.149.63"
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.112"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "'515.80"
$ws.Range("E19").Value = "  +6.45%  "
$ws.Range("D20").Value = "'6.81"
$ws.Range("E20").Value = "  +3.77%  "
$ws.Range("D21").Value = "'13.98"
$ws.Range("E21").Value = "  +3.04%  "
$ws.Range("D22").Value = "'0.714"
$ws.Range("E22").Value = "  +4.64%  "
$ws.Range("D23").Value = "'7.43"
$ws.Range("E23").Value = "  +4.21%  "
$ws.Range("D24").Value = "'12.73"
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("D25").Value = "'78.91"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'8.87"
$ws.Range("E27").Value = "  +15.49%  "
$ws.Range("D28").Value = "'2.82"
$ws.Range("E28").Value = "  +4.74%  "
$ws.Range("D29").Value = "'2.14"
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'26.54"
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("D32").Value = "'2.59"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "'1.13"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").Value = "'551.61"
$ws.Range("E34").Value = "  -5.24%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'6.07"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'5.37"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "'53.88"
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("D38").Value = "'0.0433"
$ws.Range("E38").Value = "  +6.46%  "
$ws.Range("D39").Value = "'0.0823"
$ws.Range("E39").Value = "  +4.45%  "
$ws.Range("D40").Value = "3.157.36"
$ws.Range("E40").Value = "  +8.12%  "
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Value = "'2.73"
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("D43").Value = "'8.27"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").Value = "'0.264"
$ws.Range("E44").Value = "  +9.61%  "
$ws.Range("D45").Value = "'2.22"
$ws.Range("E45").Value = "  +7.20%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'25.22"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").Value = "'120.63"
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "0.0₃0516"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("E51").Value = "  +2.97%  "
